$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.555.28'
$ws.Range("E2").Value = '  +0.34%  '

# Row 3
$ws.Range("D3").Value = '2.665.66'
$ws.Range("E3").Value = '  +0.86%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.68'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.97'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.48%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.575'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.97'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +10.45%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -1.70%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.338'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.24%  '

# Row 12
$ws.Range("E12").Value = '  +2.17%  '

# Row 13
$ws.Range("D13").Value = '3.130.62'
$ws.Range("E13").Value = '  +0.78%  '

# Row 14
$ws.Range("D14").Value = '59.610.56'
$ws.Range("E14").Value = '  +0.47%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.37'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +1.51%  '

# Row 16
$ws.Range("E16").Value = '  -0.02%  '

# Row 17
$ws.Range("D17").Value = '2.638.69'
$ws.Range("E17").Value = '  -0.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '343.20'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -1.82%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.45'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -1.29%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.47'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +1.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.38'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +2.38%  '

# Row 22
$ws.Range("E22").Value = '  +0.02%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.05'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +2.79%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.169'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +2.00%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.415'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -0.25%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.995'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -0.41%  '

# Row 27
$ws.Range("D27").Value = '0.0₃0813'
$ws.Range("E27").Value = '  +1.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.16'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +0.00%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.73'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +3.59%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -0.17%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.60'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +1.32%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.91'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.15'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.33%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.24'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +4.20%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.21'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +2.92%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.909'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -3.70%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.896'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +3.64%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.50'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +1.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.90'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +0.73%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.62'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -1.51%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.627'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +4.28%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.03'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +2.32%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +0.06%  '

# Row 44
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '276.42'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.79%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0980'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -1.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0543'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +2.34%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.93'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +4.36%  '

# Row 48
$ws.Range("D48").Value = '2.076.58'
$ws.Range("E48").Value = '  +0.14%  '

# Row 49
$ws.Range("E49").Value = '  +1.98%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.29'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +1.32%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0231'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.38%  '

